$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 242 (existing rows 242..334 shift down to 244..336).
$ws.Rows.Item(242).Insert()
$ws.Rows.Item(242).Insert()

# New record at row 242
$ws.Range("A242").Value = 10
$ws.Range("B242").Value = "Vega Modelo de Temuco"
$ws.Range("C242").Value = "La Araucanía"
$ws.Range("D242").Value = 45215
$ws.Range("E242").Value = 9
$ws.Range("F242").Value = "Fruta"
$ws.Range("G242").Value = 100101
$ws.Range("H242").Value = "Berries"
$ws.Range("I242").Value = 100112025
$ws.Range("J242").Value = "Frutilla"
$ws.Range("K242").Value = "Sin especificar"
$ws.Range("L242").Value = "Primera"
$ws.Range("M242").Value = 1400
$ws.Range("N242").Value = 10000
$ws.Range("O242").Value = 11000
$ws.Range("P242").Value = 10429
$ws.Range("Q242").Value = '$/bandeja 7 kilos'
$ws.Range("R242").Value = "Provincia de Melipilla"
$ws.Range("S242").Value = 1490
$ws.Range("T242").Value = 7

# New record at row 243
$ws.Range("A243").Value = 10
$ws.Range("B243").Value = "Vega Modelo de Temuco"
$ws.Range("C243").Value = "La Araucanía"
$ws.Range("D243").Value = 45215
$ws.Range("E243").Value = 9
$ws.Range("F243").Value = "Fruta"
$ws.Range("G243").Value = 100101
$ws.Range("H243").Value = "Berries"
$ws.Range("I243").Value = 100112025
$ws.Range("J243").Value = "Frutilla"
$ws.Range("K243").Value = "Sin especificar"
$ws.Range("L243").Value = "Segunda"
$ws.Range("M243").Value = 250
$ws.Range("N243").Value = 9000
$ws.Range("O243").Value = 9000
$ws.Range("P243").Value = 9000
$ws.Range("Q243").Value = '$/bandeja 7 kilos'
$ws.Range("R243").Value = "Provincia de Melipilla"
$ws.Range("S243").Value = 1286
$ws.Range("T243").Value = 7
